$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix stray leading space in "Binary Search" (row 79) ---
$ws.Range("B79").Value = "Binary Search"

# --- Insert a new row (994. Rotting Oranges) under "733. Flood Fill" ---
# CopyOrigin = xlFormatFromLeftOrAbove (0) so the new row inherits formatting
# from the row above it (matches the s="5"/s="14" styles used by sibling rows).
$ws.Rows(106).Insert(-4121, 0)
$ws.Range("A106").Value = 994
$ws.Range("B106").Value = "Rotting Oranges"
$ws.Range("C106").Value = "Medium"
$ws.Range("D106").Value = "Multi-Source BFS , Deque"
$ws.Range("E106").Value = 45834
$ws.Range("F106").Value = "Python"

# --- Re-point the "127. Word Ladder" hyperlink at its new row (115) ---
$ws.Range("B114").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B115"), "https://leetcode.com/problems/word-ladder/", "", "https://leetcode.com/problems/word-ladder/", "127. Word Ladder")

# --- Give F1 the same header style (s="3") as the other header cells ---
$ws.Range("D1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update the selection to match the author's final cursor position ---
$ws.Range("B88").Select()
